# work on training preview UI
#
# - row 53: hide the row and flip E53/F53 from 0 to 1
# - append two new backlog items (rows 134 and 135, both hidden) and one
#   trailing blank row (136, visible) to the bottom of the backlog table
# - grow Table1 / its AutoFilter to cover the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 53: now "in progress" for E/F, and collapsed -------------------
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = 1
$ws.Rows(53).Hidden = $true

# --- new backlog rows: start by cloning the formatting of the last -------
# --- existing row (133) onto B:F of the three new rows --------------------
$ws.Range("B133:F133").Copy()
$ws.Range("B134:F136").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A134").Value = 'افزودن بخش تنظيمات مربوط به تاريخ (تاريخ پايه اي سيستم، روز اول هفته، روزهاي تعطيل رسمي و ...)'
$ws.Range("B134").Value = 'سوم'
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 0
$ws.Range("F134").Value = 0

$ws.Range("A135").Value = 'افزودن بخش'

# Row 136 is an intentionally blank trailing row (kept visible), while rows
# 134/135 stay hidden like the rest of the backlog detail rows.
$ws.Rows(134).Hidden = $true
$ws.Rows(135).Hidden = $true

# --- grow the table + autofilter to include the new rows -----------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F136"))

# --- selection, purely cosmetic, matches where the author ended up -------
$ws.Range("B142").Select()
